# release V2021-1-2, happy new year
# Swap the "hospital" (F8766D / red) and "pharmacy" (619CFF / blue) pie
# slices' legend swatch colors and legend labels on the survey-type pie
# chart.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$grp = $s.Shapes.Item(1)

# ggplot2-style default palette colors used by the two swapped categories.
$red = 0xF8766D
$blue = 0x619CFF

# PowerPoint's RGB long is packed as 0xBBGGRR, so rebuild each color in
# that byte order from the RRGGBB hex values above.
function ToBGR($rrggbb) {
    $r = ($rrggbb -band 0xFF0000) -shr 16
    $g = ($rrggbb -band 0x00FF00) -shr 8
    $b = ($rrggbb -band 0x0000FF)
    return ($b * 0x10000) + ($g * 0x100) + $r
}

$redBGR = ToBGR $red
$blueBGR = ToBGR $blue

# Legend color swatches (small rectangles): rc58 was red, becomes blue;
# rc60 was blue, becomes red.
$rc58 = $grp.GroupItems.Item("rc58")
$rc58.Fill.ForeColor.RGB = $blueBGR
$rc58.Fill.Transparency = 0

$rc60 = $grp.GroupItems.Item("rc60")
$rc60.Fill.ForeColor.RGB = $redBGR
$rc60.Fill.Transparency = 0

# Legend text labels: tx61 was "医院" (hospital), becomes "药店" (pharmacy);
# tx63 was "药店" (pharmacy), becomes "医院" (hospital).
$tx61 = $grp.GroupItems.Item("tx61")
$tx61.TextFrame.TextRange.Text = "药店"

$tx63 = $grp.GroupItems.Item("tx63")
$tx63.TextFrame.TextRange.Text = "医院"
